# Generate Report for Handback
# ---------------------------------------------------------------
# The handback pass completed for both locales, so:
#   * the per-file "Status" flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown
#     (Overview sheet + each locale sheet's Status column),
#   * each locale sheet grows two new populated columns:
#       F "Latest Target File"   (the source .md, same file/link as A)
#       G "Latest Handback File" (the handed-back .xlf, same file/link as D)
#   * each locale sheet's "Latest Handback DateTime" (H) is stamped with
#     the real handback time (different per locale).
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276   # OLE (BGR) value of RGB(0x64,0x95,0xED) - matches the workbook's HyperLink style

function Style-AsHyperlink($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------
# Overview sheet: the Status column for both locales is driven by the
# same shared text, so flip every cell that currently reads the old
# status to the new one.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $statusNew
$wsOverview.Range("C2").Value2 = $statusNew
$wsOverview.Range("B3").Value2 = $statusNew
$wsOverview.Range("C3").Value2 = $statusNew

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value2 = $statusNew
$ws.Range("C3").Value2 = $statusNew

# New "Latest Target File" / "Latest Handback File" columns, row 2
# (97f02eb6-3868-45bf-bdc9-eec5efc9cd25 file pair)
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/efe7d3be14cdc7b13d7c4dfd3ec562751ba01af6/e2e/97f02eb6-3868-45bf-bdc9-eec5efc9cd25.md", "", "", "97f02eb6-3868-45bf-bdc9-eec5efc9cd25.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/383b5142e30563a012fc9a96748f864c45160ed8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97f02eb6-3868-45bf-bdc9-eec5efc9cd25.8a3a15b8aeeaa436431f53eb623dea5b0c7d03f7.zh-cn.xlf", "", "", "97f02eb6-3868-45bf-bdc9-eec5efc9cd25.8a3a15b8aeeaa436431f53eb623dea5b0c7d03f7.zh-cn.xlf") | Out-Null

# New columns, row 3 (ae94a271-8b44-4652-a391-beb04cb338c6 file pair)
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/efe7d3be14cdc7b13d7c4dfd3ec562751ba01af6/e2e/ae94a271-8b44-4652-a391-beb04cb338c6.md", "", "", "ae94a271-8b44-4652-a391-beb04cb338c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/383b5142e30563a012fc9a96748f864c45160ed8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ae94a271-8b44-4652-a391-beb04cb338c6.3db71f5cf1e4183c10e3cb06adb99d5c0c7b465e.zh-cn.xlf", "", "", "ae94a271-8b44-4652-a391-beb04cb338c6.3db71f5cf1e4183c10e3cb06adb99d5c0c7b465e.zh-cn.xlf") | Out-Null

Style-AsHyperlink $ws.Range("F2:G3")

# zh-cn was handed back at 17:11:20
$ws.Range("H2").Value2 = "2016-03-22 17:11:20"
$ws.Range("H3").Value2 = "2016-03-22 17:11:20"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("C2").Value2 = $statusNew
$ws2.Range("C3").Value2 = $statusNew

# New "Latest Target File" / "Latest Handback File" columns, row 2
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/efe7d3be14cdc7b13d7c4dfd3ec562751ba01af6/e2e/97f02eb6-3868-45bf-bdc9-eec5efc9cd25.md", "", "", "97f02eb6-3868-45bf-bdc9-eec5efc9cd25.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e32d984def9ea22abe02e74829c1f908473842b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97f02eb6-3868-45bf-bdc9-eec5efc9cd25.8a3a15b8aeeaa436431f53eb623dea5b0c7d03f7.de-de.xlf", "", "", "97f02eb6-3868-45bf-bdc9-eec5efc9cd25.8a3a15b8aeeaa436431f53eb623dea5b0c7d03f7.de-de.xlf") | Out-Null

# New columns, row 3
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/efe7d3be14cdc7b13d7c4dfd3ec562751ba01af6/e2e/ae94a271-8b44-4652-a391-beb04cb338c6.md", "", "", "ae94a271-8b44-4652-a391-beb04cb338c6.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e32d984def9ea22abe02e74829c1f908473842b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ae94a271-8b44-4652-a391-beb04cb338c6.3db71f5cf1e4183c10e3cb06adb99d5c0c7b465e.de-de.xlf", "", "", "ae94a271-8b44-4652-a391-beb04cb338c6.3db71f5cf1e4183c10e3cb06adb99d5c0c7b465e.de-de.xlf") | Out-Null

Style-AsHyperlink $ws2.Range("F2:G3")

# de-de was handed back at 17:11:27
$ws2.Range("H2").Value2 = "2016-03-22 17:11:27"
$ws2.Range("H3").Value2 = "2016-03-22 17:11:27"

Write-Host "Handback report generated for zh-cn and de-de."
